$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: volume/number and week-covering dates ---
$ws.Range("A8").Value = "Volume 32   Number  27"
$ws.Range("C9").Value = "Report Covering the Week  6/30/2025  Through  7/6/2025"

# --- Crime statistics table updates (rows 14-30) ---
$ws.Range("D14").Value = "0"
$ws.Range("E14").Value = "***.*"
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 2
$ws.Range("J15").Value = 5
$ws.Range("K15").Value = 60
$ws.Range("C16").Value = 6
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 28
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = 64.705882352941
$ws.Range("I16").Value = 116
$ws.Range("J16").Value = 126
$ws.Range("K16").Value = -7.936507936507
$ws.Range("L16").Value = 63.380281690140
$ws.Range("M16").Value = 13.725490196078
$ws.Range("N16").Value = -81.045751633986
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = 10
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 34
$ws.Range("H17").Value = -20.588235294117
$ws.Range("I17").Value = 185
$ws.Range("J17").Value = 171
$ws.Range("K17").Value = 8.187134502923
$ws.Range("L17").Value = 44.53125
$ws.Range("M17").Value = 208.333333333333
$ws.Range("N17").Value = -11.483253588516
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = -18.75
$ws.Range("I18").Value = 110
$ws.Range("J18").Value = 108
$ws.Range("K18").Value = 1.851851851851
$ws.Range("L18").Value = -6.779661016949
$ws.Range("M18").Value = 83.333333333333
$ws.Range("N18").Value = -73.429951690821
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -37.5
$ws.Range("F19").Value = 40
$ws.Range("G19").Value = 64
$ws.Range("H19").Value = -37.5
$ws.Range("I19").Value = 337
$ws.Range("J19").Value = 377
$ws.Range("K19").Value = -10.610079575596
$ws.Range("L19").Value = -3.714285714285
$ws.Range("M19").Value = 45.887445887445
$ws.Range("N19").Value = -35.316698656429
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 44.444444444444
$ws.Range("I20").Value = 52
$ws.Range("J20").Value = 39
$ws.Range("K20").Value = 33.333333333333
$ws.Range("L20").Value = 44.444444444444
$ws.Range("M20").Value = 116.666666666667
$ws.Range("N20").Value = -85.595567867036
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = 34
$ws.Range("E21").Value = -2.941176470588
$ws.Range("F21").Value = 121
$ws.Range("G21").Value = 143
$ws.Range("H21").Value = -15.384615384615
$ws.Range("I21").Value = 808
$ws.Range("J21").Value = 828
$ws.Range("K21").Value = -2.415458937198
$ws.Range("L21").Value = 14.447592067988
$ws.Range("M21").Value = 68.684759916492
$ws.Range("N21").Value = -62.047909816815
$ws.Range("C22").Value = 3
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = -16.666666666666
$ws.Range("I22").Value = 32
$ws.Range("J22").Value = 25
$ws.Range("K22").Value = 28
$ws.Range("L22").Value = 3.225806451612
$ws.Range("M22").Value = -5.882352941176
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 3
$ws.Range("H23").Value = 200
$ws.Range("I23").Value = 18
$ws.Range("J23").Value = 23
$ws.Range("K23").Value = -21.739130434782
$ws.Range("L23").Value = -30.769230769230
$ws.Range("M23").Value = 63.636363636363
$ws.Range("C24").Value = 44
$ws.Range("D24").Value = 48
$ws.Range("E24").Value = -8.333333333333
$ws.Range("F24").Value = 142
$ws.Range("G24").Value = 178
$ws.Range("H24").Value = -20.224719101123
$ws.Range("I24").Value = 857
$ws.Range("J24").Value = 1117
$ws.Range("K24").Value = -23.276633840644
$ws.Range("L24").Value = -17.037754114230
$ws.Range("M24").Value = 17.076502732240
$ws.Range("C25").Value = 33
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 125
$ws.Range("G25").Value = 156
$ws.Range("H25").Value = -19.871794871794
$ws.Range("I25").Value = 662
$ws.Range("J25").Value = 1006
$ws.Range("K25").Value = -34.194831013916
$ws.Range("L25").Value = -28.509719222462
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = -37.5
$ws.Range("F26").Value = 40
$ws.Range("G26").Value = 39
$ws.Range("H26").Value = 2.564102564102
$ws.Range("I26").Value = 220
$ws.Range("J26").Value = 245
$ws.Range("K26").Value = -10.204081632653
$ws.Range("L26").Value = 17.021276595744
$ws.Range("M26").Value = 18.918918918918
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 3
$ws.Range("J27").Value = 9
$ws.Range("K27").Value = 33.333333333333
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 40
$ws.Range("I28").Value = 46
$ws.Range("J28").Value = 27
$ws.Range("K28").Value = 70.370370370370
$ws.Range("L28").Value = 84
$ws.Range("C29").Value = 1
$ws.Range("F29").Value = 1
$ws.Range("I29").Value = 2
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = -66.666666666666
$ws.Range("M29").Value = -50
$ws.Range("N29").Value = -91.666666666666
$ws.Range("C30").Value = 1
$ws.Range("F30").Value = 1
$ws.Range("I30").Value = 2
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = -60
$ws.Range("M30").Value = 0
$ws.Range("N30").Value = -89.473684210526
